# Fixed no perfect maze issue
# Adds a new hour-log entry to row 10 describing the fix for the
# "no perfect maze" bug.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

$ws.Range("A10").Value = "Fixing no perfect maze option"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 44943
$ws.Range("D10").Value = 'Ater running a few test, i noticed that the algortihm did not always produce a "perfect maze". I later found out that this had to do with runnig multiple maze runners at the same time. Somehow the index of the cells in the currentCellPath did not always align. That is why i decided to only use 1 maze runner for now, but also still keep the same code structure. This way i could easily add more maze runners in the future, for when i have fixed the alignment issue.'

# Restore the view state saved with the workbook.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("D21").Select()
